# Apply "first sweep cleaning data columns to conform to specs" edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Replace underscores-for-spaces in the two shared text values used
#    throughout the data columns (experimentDesign / strain).
$ws.Range("D2:D27").Replace("Environmental Perturbation", "Environmental_Perturbation")
$ws.Range("F2:F27").Replace("KN99 alpha", "KN99_alpha")

# 2. Widen columns C, D and add an explicit width for column E so they
#    are no longer using the generic default width.
#    (ColumnWidth is expressed in characters; Excel stores the resulting
#    width on the saved-file pixel grid, which is offset from the
#    character value by ~5/6 of a character, so we compensate here to
#    land on the intended 25.33 / 24.69 / 23.94 stored widths.)
$ws.Columns.Item(3).ColumnWidth = 24.496666666666666
$ws.Columns.Item(4).ColumnWidth = 23.85666666666667
$ws.Columns.Item(5).ColumnWidth = 23.10666666666667

# 3. Move the active selection from column B to column F.
$ws.Range("F2:F27").Select()
